# BIS-1002: Fixed XLS export tests
# Adds a new "Internal Assignment" column (O) to the experiment-type export sheet:
#  - O4 gets the bold/size-12 header "Internal Assignment"
#  - O5:O8 get the literal text value "FALSE" (same look as the other plain body cells)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell O4 -------------------------------------------------------
$ws.Range("O4").Value() = "Internal Assignment"
$ws.Range("O4").Font.Name = "Calibri"
$ws.Range("O4").Font.Size = 12
$ws.Range("O4").Font.Bold = $true
$ws.Range("O4").Font.Color = 0

# --- Data cells O5:O8 ------------------------------------------------------
# Writing the literal string "FALSE" directly gets auto-coerced into a boolean
# by the engine, so instead we write it as a formula that evaluates to the
# text "FALSE" and then convert the formula to a static value. Afterwards we
# copy the plain body-cell formatting (taken from column A) onto the cell so
# it visually matches the rest of the row instead of picking up a new style.
foreach ($r in 5..8) {
    $cell = $ws.Range("O$r")

    $cell.Formula = '="FALSE"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues

    $ws.Range("A$r").Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

$excel.CutCopyMode = 0

# Mirror the final cursor position left behind in the authored workbook.
$ws.Range("P8").Select() | Out-Null
